# Generate Report for Handback
# Populates the "Latest Target File" (F) and "Latest Handback File" (G) columns
# for the zh-cn and de-de handback sheets, flips the Status to the "handed
# back" message, and stamps the handback datetime.

$wb = $excel.ActiveWorkbook

$HYPERLINK_UNDERLINE = 2        # xlUnderlineStyleSingle
$HYPERLINK_COLOR = 15570276     # RGB(100,149,237) == #6495ED, stored BGR-packed

$sheets = @(
    @{
        Name = "zh-cn"
        XlfDisplay = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        AUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0a27889de9e20db720e12b6d0fa3a4582814cee3/e2e/a.md"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c0f4b14aa3a88d45ed709512f59ee8e58c66c528/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
        HandbackDateTime = "2016-03-18 16:27:36"
    },
    @{
        Name = "de-de"
        XlfDisplay = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        AUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0a27889de9e20db720e12b6d0fa3a4582814cee3/e2e/a.md"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da740527c6286269b782431b0c43928cbbc4cff0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
        HandbackDateTime = "2016-03-18 16:27:42"
    }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Status column now reflects that the handback is in sync with en-US.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Latest Target File / Latest Handback File columns, row 2.
    $ws.Range("F2").Value = "a.md"
    $ws.Range("G2").Value = $info.XlfDisplay

    # Latest Target File / Latest Handback File columns, row 3.
    $ws.Range("F3").Value = "a.md"
    $ws.Range("G3").Value = $info.XlfDisplay

    # Latest Handback DateTime for both rows.
    $ws.Range("H2").Value = $info.HandbackDateTime
    $ws.Range("H3").Value = $info.HandbackDateTime

    # Wire up hyperlinks for the newly populated cells, in the same order
    # the handoff/target columns already use (row 2 first, then row 3).
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.AUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("G2"), $info.XlfUrl, "", "", $info.XlfDisplay)
    $ws.Hyperlinks.Add($ws.Range("F3"), $info.AUrl, "", "", "a.md")
    $ws.Hyperlinks.Add($ws.Range("G3"), $info.XlfUrl, "", "", $info.XlfDisplay)

    # Match the blue/underlined "HyperLink" look already used by the other
    # linked cells (A/B/D) in this sheet.
    $ws.Range("F2:G3").Font.Underline = $HYPERLINK_UNDERLINE
    $ws.Range("F2:G3").Font.Color = $HYPERLINK_COLOR
}

# The Overview sheet mirrors the same "Status" text for each language via a
# shared string, so keep it in sync with the handback message too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

Write-Host "Handback report generated for zh-cn and de-de."
